# feat: add 2022-Q4 data
#
# 1. Duplicate the "2022-Q3" sheet (positioned right after it becomes
#    "2022-Q4") and update its fund-size / position figures.
# 2. Insert a new "2022-Q4" row at the top of the "总计" (totals) sheet's
#    data, pushing the existing 2022-Q3 / 2022-Q2 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value that LOOKS numeric ("2.65") but must stay a text
# cell (matches the source data, which stores these figures as strings).
# ---------------------------------------------------------------------
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# =======================================================================
# Part 1: create the new "2022-Q4" sheet from a copy of "2022-Q3"
# =======================================================================
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)                 # new copy is inserted immediately before $q3
$q4 = $wb.ActiveSheet
$q4.Name = "2022-Q4"

Set-TextValue $q4 "D2" "2.65"
Set-TextValue $q4 "E2" "94.48"
Set-TextValue $q4 "F2" "4.70"
Set-TextValue $q4 "G2" "0.1246"

Set-TextValue $q4 "D3" "0.52"
Set-TextValue $q4 "E3" "94.48"
Set-TextValue $q4 "F3" "4.70"
Set-TextValue $q4 "G3" "0.0244"

# =======================================================================
# Part 2: insert the "2022-Q4" row into the "总计" totals sheet
# =======================================================================
$total = $wb.Worksheets.Item("总计")

# shift existing row 3 ("2022-Q2") down to row 4
$total.Range("A3:D3").Copy()
$total.Range("A4:D4").PasteSpecial(-4163)   # xlPasteValues
$total.Range("A3:D3").Copy()
$total.Range("A4:D4").PasteSpecial(-4122)   # xlPasteFormats

# shift existing row 2 ("2022-Q3") down to row 3
$total.Range("A2:D2").Copy()
$total.Range("A3:D3").PasteSpecial(-4163)   # xlPasteValues
$total.Range("A2:D2").Copy()
$total.Range("A3:D3").PasteSpecial(-4122)   # xlPasteFormats

# write the new "2022-Q4" row in row 2
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.15

# fix up the running index column (A) for the rows that moved down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# restore the originally-active tab (2022-Q2, the last sheet)
$wb.Worksheets.Item("2022-Q2").Activate()
